# "slides based on Geoff's talk"
#
# - The "slides" sheet is rewritten with a new set of slide rows (3-19),
#   replacing the old slide list.
# - A new sheet "old_slides" is appended (after "values") that preserves
#   the original slide id list (what used to be in "slides" A3:A17).
# - The "values" sheet is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("slides")

# --- capture the existing slide ids (A3:A17) before we overwrite them,
#     so we can stash them on the new "old_slides" sheet. ---
$oldIds = New-Object System.Collections.ArrayList
for ($r = 3; $r -le 17; $r++) {
    [void]$oldIds.Add($ws.Cells.Item($r, 1).Value2)
}

# --- clear out the old slide rows (A3:I17) and write the new data ---
$ws.Range("A3:I19").ClearContents()

$ws.Range("A3").Value = "in-this-together"
$ws.Range("B3").Value = "slide-dark"
$ws.Range("D3").Value = "in-this-together.jpg"

$ws.Range("A4").Value = "my-kind-of-town"
$ws.Range("A5").Value = "broken-fields"
$ws.Range("A6").Value = "multiple-values"
$ws.Range("A7").Value = "different-references"
$ws.Range("A8").Value = "encoded-values"
$ws.Range("A9").Value = "what-to-do"
$ws.Range("A10").Value = "rule1"
$ws.Range("A11").Value = "start-clean"
$ws.Range("A12").Value = "toolkit"
$ws.Range("A13").Value = "homework"
$ws.Range("A14").Value = "inspection"
$ws.Range("A15").Value = "quick-stats1"

$ws.Range("A16").Value = "quick-stats2"
$ws.Range("D16").Value = "refine_facet.png"

$ws.Range("A17").Value = "quick-stats3"
$ws.Range("D17").Value = "refine_cluster.png"

$ws.Range("A18").Value = "decompose"

$ws.Range("A19").Value = "do-i-have-to"
$ws.Range("D19").Value = "do-i-have-to.jpg"

# --- freeze panes / selection roughly follow the new bottom of the sheet ---
$ws.Application.Goto($ws.Range("A1"), $true)
$ws.Range("B10").Select()
$ws.Application.ActiveWindow.FreezePanes = $true
$ws.Range("B20").Select()

# --- add the "old_slides" sheet after "values", holding the previous
#     slide ids that used to live on "slides" ---
$valuesSheet = $wb.Worksheets.Item("values")
$oldSheet = $wb.Worksheets.Add($null, $valuesSheet)
$oldSheet.Name = "old_slides"

for ($i = 0; $i -lt $oldIds.Count; $i++) {
    $oldSheet.Cells.Item($i + 1, 1).Value = $oldIds[$i]
}

$ws.Select()
